# Add a second column (B) of values next to the existing column A data,
# and move the active cell/selection from B6 to C4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 2
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 2

$ws.Range("C4").Select()
